# daily auto push: 2026-02-19 19:05 UTC
#
# Insert one new log row at row 833 (pushing the existing rows 833-874 down
# to 834-875) and populate it with the new day's first reading:
#   2026/02/20, 金, 0, 201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 833:874 down by one to make room for the new entry.
$ws.Rows(833).Insert()

# Column A holds dates as literal text (e.g. "2026/12/29"), not real date
# serials. Force text formatting before assigning so Excel's automatic
# date-recognition doesn't convert the string into a date value, then
# restore the default "Normal" style so the new cell matches its siblings
# (which carry no explicit style).
$ws.Range("A833").NumberFormat = "@"
$ws.Range("A833").Value = "2026/02/20"
$ws.Range("A833").Style = "Normal"

$ws.Range("B833").Value = "金"
$ws.Range("C833").Value = 0
$ws.Range("D833").Value = 201
